# Append the new profit data row (row 33) for the 2025-12-27 run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 33

# Column A holds the date as plain text (matching the rest of the sheet),
# so force a text format before assigning, then restore the default style
# so no stray number format sticks to the cell.
$cellA = $ws.Range("A" + $row)
$cellA.NumberFormat = "@"
$cellA.Value = "12/27/2025"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 12167.42
$ws.Cells.Item($row, 3).Value = 0.206830545967442
$ws.Cells.Item($row, 4).Value = 0.793169454032558
$ws.Cells.Item($row, 5).Value = -137.02
$ws.Cells.Item($row, 6).Value = -25.61
$ws.Cells.Item($row, 7).Value = -20940.93
$ws.Cells.Item($row, 8).Value = -68.45
$ws.Cells.Item($row, 9).Value = -486.26
$ws.Cells.Item($row, 10).Value = -16.19
